$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift data rows 2-10 down into rows 3-11 (each row takes on the values
# previously held by the row above it), then write the new values for
# row 2 (the newest quarter of data) at the top. Column A (period labels)
# is left untouched.

for ($r = 10; $r -ge 2; $r--) {
    $src = $r
    $dst = $r + 1
    $srcVal = $ws.Range("B$src`:G$src").Value()
    $ws.Range("B$dst`:G$dst").Value = $srcVal
}

$ws.Range("B2").Value = -0.02907897629796788
$ws.Range("C2").Value = 0.3131278957257717
$ws.Range("D2").Value = 0.181524606355785
$ws.Range("E2").Value = 0.4260570458938391
$ws.Range("F2").Value = 0.43998257208981
$ws.Range("G2").Value = 15
